$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26..171 down to 27..172
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new weekly data point
$ws.Range("A26").Value = 5
$ws.Range("B26").Value = "Macroferia Regional de Talca"
$ws.Range("C26").Value = "Maule"
$ws.Range("D26").Value = 44462
$ws.Range("E26").Value = 7
$ws.Range("F26").Value = 100112009
$ws.Range("G26").Value = "Acelga"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 600
$ws.Range("K26").Value = 2000
$ws.Range("L26").Value = 2000
$ws.Range("M26").Value = 2000
$ws.Range("N26").Value = "$/docena de atados (4 kilos)"
$ws.Range("O26").Value = "Región del Maule"
$ws.Range("P26").Value = 500
$ws.Range("Q26").Value = 4
$ws.Range("R26").Value = "Hortaliza"
